# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# The "municipio-nombre" column (M) metadata block is updated so that it is
# curated/treated as a dimension, matching the pattern already used by the
# neighbouring "provincia-nombre" (N) and "comarca-nombre" (O) columns:
#   M2: iaest-measure:municipio-nombre -> sdmx-dimension:refArea
#   M3: medida                         -> dim
#   M4: xsd:int                        -> URI-Municipio

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = "sdmx-dimension:refArea"
$ws.Range("M3").Value = "dim"
$ws.Range("M4").Value = "URI-Municipio"
